$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "29.874.03"
$ws.Range("E2").Value = "  -1.21%  "

# Row 3
$ws.Range("D3").Value = "1.894.19"
$ws.Range("E3").Value = "  -1.23%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.001"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  -0.08%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.7755"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -4.11%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "244.75"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +0.19%  "

# Row 7
$ws.Range("E7").Value = "  -0.07%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3148"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  -2.89%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07585"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +4.77%  "

# Row 10
$ws.Range("E10").Value = "  -5.45%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.08106"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +0.17%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.7731"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  -1.95%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.493"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  +1.54%  "

# Row 14
$ws.Range("D14").Value = "1.806.35"
$ws.Range("E14").Value = "  -6.13%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "92.43"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -1.56%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "6.251"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  +3.20%  "

# Row 17
$ws.Range("D17").Value = "29.843.25"
$ws.Range("E17").Value = "  -1.28%  "

# Row 18
$ws.Range("E18").Value = "  -1.35%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000007952"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +1.35%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "244.43"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -2.09%  "

# Row 21
$ws.Range("B21").Value = "Dai"
$ws.Range("C21").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.9998"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -0.09%  "

# Row 22
$ws.Range("B22").Value = "Chainlink"
$ws.Range("C22").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "8.130"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -1.23%  "

# Row 23
$ws.Range("D23").Value = "2.126.01"
$ws.Range("E23").Value = "  -1.15%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.9999"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -0.19%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.1567"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -5.75%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.468"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -0.20%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "162.83"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -3.07%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.81"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -1.05%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.048"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -5.10%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.444"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +4.16%  "

# Row 31
$ws.Range("E31").Value = "  +0.03%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.498"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +3.88%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.100"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -1.12%  "

# Row 34
$ws.Range("E34").Value = "  -3.99%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.262"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -2.49%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7602"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +1.40%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.001"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -0.21%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.644"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -3.22%  "

# Row 39
$ws.Range("E39").Value = "  -1.55%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.788"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -1.12%  "

# Row 41
$ws.Range("D41").Value = "1.164.77"
$ws.Range("E41").Value = "  +14.21%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "73.96"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -0.66%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.4444"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -2.06%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "5.951"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -0.09%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.8480"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -0.38%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.0000"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -0.09%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.896"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -1.55%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.129"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +0.64%  "

# Row 49
$ws.Range("B49").Value = "Quant"
$ws.Range("C49").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "102.29"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -1.11%  "

# Row 50
$ws.Range("B50").Value = "EnergySwap"
$ws.Range("C50").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "9.978"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -0.48%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.547"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -0.97%  "

